{"js": "// \"terms and definitions update\"\n// The worksheet intro paragraph reads \"...You must have at least 50 different\n// terms/definitions by the end of the course.\" Update the minimum count from\n// 50 to 40.\nconst body = context.document.body;\nconst results = body.search(\"50\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find text \"50\" to replace.');\n}\n\n// Replace the matched \"50\" run's text with \"40\", preserving its formatting.\nresults.items[0].insertText(\"40\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# \"terms and definitions update\"\n# The worksheet intro paragraph reads \"...You must have at least 50 different\n# terms/definitions by the end of the course.\" Update the minimum count from\n# 50 to 40.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"50\"\n$find.Replacement.Text = \"40\"\n$find.Forward = $true\n$find.Wrap = 0\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.MatchWildcards = $false\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)\n"}
